# Applies the scheduled-runner market-data refresh to the profit-calculation
# columns (H:N) across every class tab, per the authoritative OOXML diff for
# Sheets/Marilith_Profits.xlsx. Only currentAveragePrice*/LevePrice*/LeveProfit*
# cells move; no formulas exist in this workbook, so each touched cell is
# written as a literal value (ClearContents for the handful of cells the diff
# removes outright).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 41.88889
$ws.Range("J9").Value = 101
$ws.Range("L9").Value = 101
$ws.Range("N9").Value = -439
$ws.Range("H15").Value = 358.42856
$ws.Range("I15").Value = 358.42856
$ws.Range("K15").Value = 1075.28568
$ws.Range("M15").Value = -906.28568
$ws.Range("H33").Value = 380
$ws.Range("I33").Value = 212.55556
$ws.Range("J33").Value = 568.375
$ws.Range("K33").Value = 212.55556
$ws.Range("L33").Value = 568.375
$ws.Range("M33").Value = 16.44443999999999
$ws.Range("N33").Value = -1026.375
$ws.Range("H58").Value = 742.44446
$ws.Range("I58").Value = 460.25
$ws.Range("K58").Value = 1380.75
$ws.Range("M58").Value = -1230.75
$ws.Range("H80").Value = 668.8
$ws.Range("I80").Value = 728.9167
$ws.Range("K80").Value = 2186.7501
$ws.Range("M80").Value = -1188.7501
$ws.Range("H83").Value = 668.8
$ws.Range("I83").Value = 728.9167
$ws.Range("K83").Value = 6560.2503
$ws.Range("M83").Value = -1568.2503
$ws.Range("H129").Value = 1943.4375
$ws.Range("I129").Value = 733
$ws.Range("J129").Value = 3011.4707
$ws.Range("K129").Value = 2199
$ws.Range("L129").Value = 9034.4121
$ws.Range("M129").Value = 2801
$ws.Range("N129").Value = -19034.4121
$ws.Range("I131").Value = 366.66666
$ws.Range("J131").Value = 250
$ws.Range("K131").Value = 1099.99998
$ws.Range("L131").Value = 750
$ws.Range("M131").Value = 3940.00002
$ws.Range("N131").Value = -10830
$ws.Range("H132").Value = 4367.143
$ws.Range("I132").Value = 2946
$ws.Range("J132").Value = 12894
$ws.Range("K132").Value = 8838
$ws.Range("L132").Value = 38682
$ws.Range("M132").Value = -6308
$ws.Range("N132").Value = -43742

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3697.8572
$ws.Range("I2").Value = 3357
$ws.Range("J2").Value = 4550
$ws.Range("K2").Value = 3357
$ws.Range("L2").Value = 4550
$ws.Range("M2").Value = -3244
$ws.Range("N2").Value = -4776
$ws.Range("H74").Value = 2837.3333
$ws.Range("I74").Value = 2837.3333
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 2837.3333
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -1963.3333
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 2837.3333
$ws.Range("I77").Value = 2837.3333
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 14186.6665
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -9818.666499999999
$ws.Range("N77").ClearContents()
$ws.Range("H110").Value = 2795.8
$ws.Range("I110").Value = 536.25
$ws.Range("J110").Value = 4302.1665
$ws.Range("K110").Value = 536.25
$ws.Range("L110").Value = 4302.1665
$ws.Range("M110").Value = 1508.75
$ws.Range("N110").Value = -8392.166499999999
$ws.Range("H116").Value = 3697.8572
$ws.Range("I116").Value = 3357
$ws.Range("J116").Value = 4550
$ws.Range("K116").Value = 3357
$ws.Range("L116").Value = 4550
$ws.Range("M116").Value = -1063
$ws.Range("N116").Value = -9138
$ws.Range("H132").Value = 842.4167
$ws.Range("I132").Value = 842.4167
$ws.Range("K132").Value = 2527.2501
$ws.Range("M132").Value = 2.749899999999798

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3697.8572
$ws.Range("I3").Value = 3357
$ws.Range("J3").Value = 4550
$ws.Range("K3").Value = 3357
$ws.Range("L3").Value = 4550
$ws.Range("M3").Value = -3243
$ws.Range("N3").Value = -4778
$ws.Range("H86").Value = 2366.5
$ws.Range("I86").Value = 2193.7273
$ws.Range("K86").Value = 2193.7273
$ws.Range("M86").Value = -1070.7273
$ws.Range("H89").Value = 2366.5
$ws.Range("I89").Value = 2193.7273
$ws.Range("K89").Value = 10968.6365
$ws.Range("M89").Value = -5352.636500000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 70.42856999999999
$ws.Range("I7").Value = 63.333332
$ws.Range("K7").Value = 63.333332
$ws.Range("M7").Value = 49.666668
$ws.Range("H69").Value = 1111
$ws.Range("I69").Value = 1111
$ws.Range("K69").Value = 1111
$ws.Range("M69").Value = -362
$ws.Range("H72").Value = 1111
$ws.Range("I72").Value = 1111
$ws.Range("K72").Value = 3333
$ws.Range("M72").Value = 411
$ws.Range("H86").Value = 4084.3333
$ws.Range("I86").Value = 3168.6667
$ws.Range("J86").Value = 5000
$ws.Range("K86").Value = 3168.6667
$ws.Range("L86").Value = 5000
$ws.Range("M86").Value = -2045.6667
$ws.Range("N86").Value = -7246
$ws.Range("H89").Value = 4084.3333
$ws.Range("I89").Value = 3168.6667
$ws.Range("J89").Value = 5000
$ws.Range("K89").Value = 15843.3335
$ws.Range("L89").Value = 25000
$ws.Range("M89").Value = -10227.3335
$ws.Range("N89").Value = -36232
$ws.Range("H107").Value = 1189.8572
$ws.Range("I107").Value = 1265
$ws.Range("J107").Value = 739
$ws.Range("K107").Value = 1265
$ws.Range("L107").Value = 739
$ws.Range("M107").Value = 655
$ws.Range("N107").Value = -4579
$ws.Range("H132").Value = 3550
$ws.Range("J132").Value = 3600
$ws.Range("L132").Value = 10800
$ws.Range("N132").Value = -15860

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1183.4667
$ws.Range("I4").Value = 1131.8518
$ws.Range("J4").Value = 1260.8889
$ws.Range("K4").Value = 3395.5554
$ws.Range("L4").Value = 3782.6667
$ws.Range("M4").Value = -3283.5554
$ws.Range("N4").Value = -4006.6667
$ws.Range("H22").Value = 3980
$ws.Range("I22").Value = 1500
$ws.Range("K22").Value = 4500
$ws.Range("M22").Value = -4331
$ws.Range("H27").Value = 3980
$ws.Range("I27").Value = 1500
$ws.Range("K27").Value = 4500
$ws.Range("M27").Value = -4398
$ws.Range("H44").Value = 1331.3334
$ws.Range("I44").Value = 1397.6
$ws.Range("K44").Value = 4192.799999999999
$ws.Range("M44").Value = -3794.799999999999
$ws.Range("H124").Value = 2500
$ws.Range("I124").Value = 2500
$ws.Range("K124").Value = 7500
$ws.Range("M124").Value = -2590

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 16859.166
$ws.Range("J57").Value = 25250
$ws.Range("L57").Value = 25250
$ws.Range("N57").Value = -26890
$ws.Range("H80").Value = 3750
$ws.Range("J80").Value = 3750
$ws.Range("L80").Value = 3750
$ws.Range("N80").Value = -5746
$ws.Range("H83").Value = 3750
$ws.Range("J83").Value = 3750
$ws.Range("L83").Value = 18750
$ws.Range("N83").Value = -28734
$ws.Range("H132").Value = 2718.125
$ws.Range("I132").Value = 1957.5
$ws.Range("K132").Value = 5872.5
$ws.Range("M132").Value = -3342.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2181.9375
$ws.Range("I46").Value = 2221.6667
$ws.Range("K46").Value = 2221.6667
$ws.Range("M46").Value = -2033.6667
$ws.Range("H51").Value = 38417.332
$ws.Range("J51").Value = 38417.332
$ws.Range("L51").Value = 38417.332
$ws.Range("N51").Value = -39373.332
$ws.Range("H55").Value = 239.28572
$ws.Range("I55").Value = 155
$ws.Range("J55").Value = 450
$ws.Range("K55").Value = 155
$ws.Range("L55").Value = 450
$ws.Range("M55").Value = 18
$ws.Range("N55").Value = -796
$ws.Range("H61").Value = 7427.6665
$ws.Range("I61").Value = 7427.6665
$ws.Range("K61").Value = 7427.6665
$ws.Range("M61").Value = -7225.6665
$ws.Range("H113").Value = 7427.6665
$ws.Range("I113").Value = 7427.6665
$ws.Range("K113").Value = 7427.6665
$ws.Range("M113").Value = -5257.6665
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 28125
$ws.Range("J54").Value = 32500
$ws.Range("L54").Value = 32500
$ws.Range("N54").Value = -33540
$ws.Range("H96").Value = 900
$ws.Range("I96").Value = 900
$ws.Range("K96").Value = 900
$ws.Range("M96").Value = 473
$ws.Range("H107").Value = 275.66666
$ws.Range("I107").Value = 275.66666
$ws.Range("K107").Value = 826.9999799999999
$ws.Range("M107").Value = 1093.00002

